$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.030.63"
$ws.Range("D3").Value = "1.832.13"
$ws.Range("E3").Value = "  -1.32%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6704"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.37%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07409"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.56%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2952"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.71"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.99%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07643"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.59%  "
$ws.Range("D12").Value = "1.832.27"
$ws.Range("E12").Value = "  -1.39%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.004"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6723"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "85.76"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.142"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.63%  "
$ws.Range("D17").Value = "29.037.54"
$ws.Range("E17").Value = "  -1.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008230"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "227.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9994"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.305"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "160.46"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.60%  "
$ws.Range("E25").Value = "  -4.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.665"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.93"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.501"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.54%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.224"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.60%  "
$ws.Range("E30").Value = "  -1.29%  "
$ws.Range("E31").Value = "  -0.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05364"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7515"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.851"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.01%  "
$ws.Range("E35").Value = "  -2.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.681"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.14%  "
$ws.Range("D37").Value = "1.288.98"
$ws.Range("E37").Value = "  -3.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01807"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.39%  "
$ws.Range("E39").Value = "  -0.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9215"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.004"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "104.26"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9982"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.08011"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +21.13%  "
$ws.Range("D45").Value = "1.978.05"
$ws.Range("E45").Value = "  -1.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5177"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.78%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000121"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.80%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.355"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.35%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "63.46"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.07%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.751"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.29%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05919"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.07%  "
